# Update cryptos list with latest price/volume data (GitHub Actions scrape refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "53.480.61"
$ws.Range("E2").Value = "  -4.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.192.19"
$ws.Range("E3").Value = "  -7.27%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "485.14"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "124.93"
$ws.Range("E6").Value = "  -3.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.520"
$ws.Range("E8").Value = "  -4.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.210.21"
$ws.Range("E9").Value = "  -6.60%  "
$ws.Range("E10").Value = "  -6.85%  "
$ws.Range("E11").Value = "  -1.44%  "
$ws.Range("E12").Value = "  -4.32%  "
$ws.Range("E13").Value = "  -3.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.582.24"
$ws.Range("E14").Value = "  -7.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.08"
$ws.Range("E15").Value = "  -1.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "53.404.77"
$ws.Range("E16").Value = "  -4.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000127"
$ws.Range("E17").Value = "  -3.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.211.86"
$ws.Range("E18").Value = "  -6.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.54"
$ws.Range("E19").Value = "  -4.83%  "
$ws.Range("E20").Value = "  -1.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "293.62"
$ws.Range("E21").Value = "  -4.64%  "
$ws.Range("E22").Value = "  -3.40%  "
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.43"
$ws.Range("E24").Value = "  -5.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.996"
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.363"
$ws.Range("E26").Value = "  -2.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.145"
$ws.Range("E27").Value = "  -1.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.294.38"
$ws.Range("E28").Value = "  -7.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.98"
$ws.Range("E29").Value = "  -3.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "165.41"
$ws.Range("E30").Value = "  -4.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.57"
$ws.Range("E31").Value = "  -4.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.998"
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.994"
$ws.Range("E33").Value = "  -0.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0660"
$ws.Range("E34").Value = "  -7.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.67"
$ws.Range("E35").Value = "  -2.08%  "
$ws.Range("E36").Value = "  -2.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.24"
$ws.Range("E37").Value = "  -2.23%  "
$ws.Range("E38").Value = "  -2.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.822"
$ws.Range("E39").Value = "  +2.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.77"
$ws.Range("E40").Value = "  -1.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.52"
$ws.Range("E41").Value = "  -5.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.364"
$ws.Range("E42").Value = "  -1.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.35"
$ws.Range("E43").Value = "  -2.50%  "
$ws.Range("E44").Value = "  -2.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "124.55"
$ws.Range("E45").Value = "  -3.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.74"
$ws.Range("E46").Value = "  +1.06%  "
$ws.Range("E47").Value = "  -2.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.531"
$ws.Range("E48").Value = "  -5.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "230.64"
$ws.Range("E49").Value = "  -3.32%  "
$ws.Range("E50").Value = "  -2.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0199"
$ws.Range("E51").Value = "  -3.44%  "

Write-Output "Applied cryptos update"
